# Fix Higiene screen bug
#
# Row 754 on the "Translation" sheet ("ABRIR PORTA / PARA COMEÇAR" door-open
# text, SingleUseId4100) was out of order - it needs to move down to just
# before the first fully-blank row (originally row 771), and every row
# between its old and new position shifts up by one to close the gap.
#
# Implemented as: remember row 754's values, delete row 754 (shifts 755..770
# up to 754..769), insert a fresh blank row at 770 (shifts things back down,
# recreating the gap in the right place), then write the remembered values
# into that new row 770.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$sourceRow = 754
$destRow = 770
$firstCol = 2   # column B
$lastCol = 9    # column I

# Remember the row's contents before we move anything.
$savedValues = @()
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $savedValues += ,$ws.Cells.Item($sourceRow, $c).Value2
}

# Remove the row, shifting everything below it up by one.
$ws.Rows.Item($sourceRow).Delete()

# Re-open a blank row just before the destination, shifting rows back down.
$ws.Rows.Item($destRow).Insert()

# Write the remembered row into its new home.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item($destRow, $c).Value2 = $savedValues[$c - $firstCol]
}
